# DeleteCustomFilter.xlsx test-data update
# Adds new "Drop 3" filter-name assertions (with slightly-mutated labels)
# to the B12 cells of each sheet, fills in the previously-empty C3 cells
# with mock numeric-looking filter ids, tweaks a couple of row heights,
# and moves the active selection/tab around (ending on DelWorklistOverview).

$wb = $excel.ActiveWorkbook

$sheetConcl     = $wb.Worksheets.Item("DelConclToApprove")
$sheetDonation  = $wb.Worksheets.Item("DelDonationInfo")
$sheetTest      = $wb.Worksheets.Item("DelTestInfo")
$sheetWlDetail  = $wb.Worksheets.Item("DelWorklistDetail")
$sheetWlOverview= $wb.Worksheets.Item("DelWorklistOverview")

# --- New B12 "filterName" values (introduces 5 new shared strings, in this order) ---
$sheetConcl.Range("B12").Value      = "ConclusionsToApprovert"
$sheetDonation.Range("B12").Value   = "DonationInformationry"
$sheetTest.Range("B12").Value       = "TestInformationrty"
$sheetWlDetail.Range("B12").Value   = "WorklistConclusionsyrt"
$sheetWlOverview.Range("B12").Value = "WorklistInformationty"

# --- Fill in the previously-empty C3 cells ---
$sheetConcl.Range("C3").Value      = "124$"
$sheetDonation.Range("C3").Value   = "123$"
$sheetTest.Range("C3").Value       = "123$"
$sheetWlDetail.Range("C3").Value   = "123$"
$sheetWlOverview.Range("C3").Value = "123$"

# --- Row height tweaks on DelConclToApprove / DelDonationInfo ---
$sheetConcl.Rows.Item(11).RowHeight    = 30
$sheetConcl.Rows.Item(12).RowHeight    = 60
$sheetDonation.Rows.Item(12).RowHeight = 45

# --- Selection / active-cell + active-sheet/tab moves ---
# (Selecting on the last sheet last makes it the active/tabSelected sheet,
# and updates the workbook's bookViews/activeTab accordingly.)
[void]$sheetConcl.Range("E18").Select()
[void]$sheetDonation.Range("F12").Select()
[void]$sheetTest.Range("H21").Select()
[void]$sheetWlDetail.Range("K16").Select()
[void]$sheetWlOverview.Range("L17").Select()
